# Auto-generated edit script applying the RSA paper expansion diff.
$word.Options.StoreRSIDOnSave = $false

function Find-ParaIndex($doc, $snippet) {
    $i = 1
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.Contains($snippet)) {
            return $i
        }
        $i = $i + 1
    }
    return -1
}

$d = $word.ActiveDocument

# 1) Background and Motivation paragraph: extend sentence about algorithmic problems
$idx = Find-ParaIndex $d 'During this course, I became fascinated'
$d.Paragraphs($idx).Range.Text = 'During this course, I became fascinated by how Large Language Models (LLMs) can assist in solving programming challenges. As someone who has participated in ACM-ICPC competitions, I often struggled with complex algorithmic problems that require deep mathematical reasoning and careful implementation. This inspired me to explore whether LLMs, combined with the Recursive Self-Aggregation (RSA) technique, could help solve these difficult problems more effectively.'

# 2) Insert new paragraph right after it (inherits spacing/justify formatting)
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$motivIdx = $idx + 1
$d.Paragraphs($motivIdx).Range.Text = 'The motivation behind this project stems from observing that competitive programming problems often have multiple layers of complexity. A single prompt to an LLM rarely produces correct solutions for hard problems. I wanted to understand how structured approaches like RSA could bridge this gap by enabling iterative problem-solving.'

# 3) 'Through my exploration' paragraph rewrite (contains an apostrophe; Range.Text avoids smart-quote autocorrect)
$idx = Find-ParaIndex $d 'Through my exploration'
$d.Paragraphs($idx).Range.Text = 'Through my exploration, I discovered that simply asking an LLM to solve a hard algorithm problem usually fails due to the model''s tendency to make logical errors or misunderstand problem constraints. The key insight from papers like Wei et al. (2022) on chain-of-thought reasoning is that breaking problems into smaller steps helps significantly. The RSA approach takes this further by recursively decomposing problems and letting the model refine its answers through multiple iterations with feedback.'

# 4) 'I learned that LLMs work best' paragraph rewrite
$idx = Find-ParaIndex $d 'I learned that LLMs work best'
$d.Paragraphs($idx).Range.Text = 'I learned that LLMs work best when given clear structure and feedback loops. The iterative nature of RSA allows the model to learn from its mistakes by analyzing test case failures. This self-correcting mechanism is crucial for handling edge cases common in competitive programming.'

# 5) 'I developed a Python program' paragraph rewrite
$idx = Find-ParaIndex $d 'I developed a Python program'
$d.Paragraphs($idx).Range.Text = 'I developed a Python program that implements the RSA workflow using Google Gemini API. The implementation follows a systematic approach with four distinct phases:'

# 6) 'Problem Analysis:' bullet rewrite
$idx = Find-ParaIndex $d 'Problem Analysis:'
$d.Paragraphs($idx).Range.Text = 'Problem Analysis: The LLM analyzes the problem statement, identifies the problem type (e.g., dynamic programming, graph theory), and breaks it into manageable subproblems.'

# 7) 'Solution Generation:' bullet rewrite
$idx = Find-ParaIndex $d 'Solution Generation:'
$d.Paragraphs($idx).Range.Text = 'Solution Generation: Based on the analysis, the LLM generates Python code with proper input/output handling and algorithm implementation.'

# 8) 'Iterative Refinement:' bullet becomes 'Testing and Validation:' (text swapped in place)
$idx = Find-ParaIndex $d 'Iterative Refinement: If tests fail'
$d.Paragraphs($idx).Range.Text = 'Testing and Validation: The generated solution is automatically tested against provided test cases with timeout handling.'

# 9) Insert new ListBullet paragraph after it carrying the new Iterative Refinement text
#    (inherits ListBullet style + spacing automatically)
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$refineIdx = $idx + 1
$d.Paragraphs($refineIdx).Range.Text = 'Iterative Refinement: If tests fail, the LLM receives detailed error information and generates an improved solution until all tests pass or maximum iterations are reached.'

# 10) 'Testing my implementation' paragraph rewrite
$idx = Find-ParaIndex $d 'Testing my implementation on several ACM-ICPC problems'
$d.Paragraphs($idx).Range.Text = 'Testing my implementation on several ACM-ICPC problems, I found that the RSA approach improved success rates significantly compared to direct prompting. Problems involving dynamic programming and graph algorithms benefited most from the recursive refinement process. The iterative feedback mechanism allowed the model to correct off-by-one errors and boundary condition issues that are common failure points.'

# 11) Insert new paragraph after it for the 'deepened my understanding' text
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$deepIdx = $idx + 1
$d.Paragraphs($deepIdx).Range.Text = 'This project deepened my understanding of both competitive programming strategies and AI-assisted problem-solving. I realized that successful AI-assisted programming lies in designing robust feedback systems that enable continuous improvement, rather than expecting perfect first-attempt solutions.'

# 12) Insert paragraph for 'Program Framework' heading (promoted to Heading2 at the end,
#     after its body-formatted successors are already created, so they don't inherit Heading2)
$d.Paragraphs($deepIdx).Range.InsertParagraphAfter()
$headIdx = $deepIdx + 1
$d.Paragraphs($headIdx).Range.Text = 'Program Framework'

# 13) Insert paragraph after the heading for the 'Below is the general structure...' sentence
$d.Paragraphs($headIdx).Range.InsertParagraphAfter()
$belowIdx = $headIdx + 1
$d.Paragraphs($belowIdx).Range.Text = 'Below is the general structure of my RSA algorithm implementation:'

# 14) Insert paragraph after that for the code block
$d.Paragraphs($belowIdx).Range.InsertParagraphAfter()
$codeIdx = $belowIdx + 1
$codeText = @'
class RSAAlgorithmSolver:
    def __init__(self, max_iterations=5):
        self.model = genai.GenerativeModel('gemini-1.5-pro')
        self.max_iterations = max_iterations
    
    def analyze_problem(self, problem_description):
        # Step 1: Use LLM to decompose problem
        # Returns: problem type, subproblems, approach
        
    def generate_solution(self, problem_analysis):
        # Step 2: Generate Python code solution
        # Returns: executable Python code
        
    def test_solution(self, code, test_cases):
        # Step 3: Run code against test cases
        # Returns: pass/fail results with details
        
    def refine_solution(self, code, test_results, analysis):
        # Step 4: RSA core - iteratively improve
        # Returns: refined code based on failures
        
    def solve(self, problem, test_cases):
        # Main workflow orchestrator
        analysis = self.analyze_problem(problem)
        code = self.generate_solution(analysis)
        for iteration in range(self.max_iterations):
            results = self.test_solution(code, test_cases)
            if all_passed(results):
                return code
            code = self.refine_solution(code, results, analysis)
        return code
'@
$lb = [char]11
$codeLines = $codeText -split "`n"
$joined = $lb + ($codeLines -join $lb) + $lb
$d.Paragraphs($codeIdx).Range.Text = $joined
$d.Paragraphs($codeIdx).Format.LineSpacingRule = 0
$d.Paragraphs($codeIdx).Alignment = 0
$codeTextRange = $d.Range($d.Paragraphs($codeIdx).Range.Start, $d.Paragraphs($codeIdx).Range.End - 1)
$codeTextRange.Font.Name = 'Courier New'
$codeTextRange.Font.Size = 9

# 15) Now promote the 'Program Framework' paragraph to Heading2 style
$d.Paragraphs($headIdx).Style = 'Heading2'

Write-Output 'edit applied successfully'
